$wb = $excel.ActiveWorkbook

# ----- Sheet: Overview -----
$ws = $wb.Worksheets.Item("Overview")

# Update cell values that changed
$ws.Range("A2").Value = "ffff442e56c6-4387-47f3-a4e8-312386de752c.md"
$ws.Range("A3").Value = "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md"
$ws.Range("A4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

# Rebuild hyperlinks: same target URLs, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md", "", "", "ffff442e56c6-4387-47f3-a4e8-312386de752c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffff442e56c6-4387-47f3-a4e8-312386de752c.md", "", "", "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffffff60694c53-22b2-4af5-9054-0144e0a92462.md", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/.localization-config", "", "", ".localization-config") | Out-Null

# ----- Sheet: zh-cn -----
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values that changed
$ws.Range("A2").Value = "ffff442e56c6-4387-47f3-a4e8-312386de752c.md"
$ws.Range("C2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-08 21:19:28"
$ws.Range("E2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.md"
$ws.Range("F2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-08 21:20:16"
$ws.Range("A3").Value = "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md"
$ws.Range("A4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-08 21:25:38"
$ws.Range("E4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md"
$ws.Range("F4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-08 21:24:42"

# Rebuild hyperlinks: same target URLs, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md", "", "", "ffff442e56c6-4387-47f3-a4e8-312386de752c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/546c903f24b37b7ce8d305f96aa80bc44ec16bf1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/07d1fed9bd0af1e4e507aa4192a687fc3298ec94/e2e/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aa4aae842e66e141e2c72d335e99215d668e88d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffff442e56c6-4387-47f3-a4e8-312386de752c.md", "", "", "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a416264f3888d2274630f0d2717e425b54cb7dec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f52d43d8f7f461afeef3e0d68b9f0ba41e408e98/e2e/314fba8a-d786-4633-a6c9-f710a8d5fa37.md", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9859adee7f9c6a443259957de0ffc59badba5dbb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffffff60694c53-22b2-4af5-9054-0144e0a92462.md", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a416264f3888d2274630f0d2717e425b54cb7dec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f52d43d8f7f461afeef3e0d68b9f0ba41e408e98/e2e/314fba8a-d786-4633-a6c9-f710a8d5fa37.md", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9859adee7f9c6a443259957de0ffc59badba5dbb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.zh-cn.xlf", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/.localization-config", "", "", ".localization-config") | Out-Null

# ----- Sheet: de-de -----
$ws = $wb.Worksheets.Item("de-de")

# Update cell values that changed
$ws.Range("A2").Value = "ffff442e56c6-4387-47f3-a4e8-312386de752c.md"
$ws.Range("C2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf"
$ws.Range("D2").Value = "2016-03-08 21:19:35"
$ws.Range("E2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.md"
$ws.Range("F2").Value = "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf"
$ws.Range("G2").Value = "2016-03-08 21:20:38"
$ws.Range("A3").Value = "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md"
$ws.Range("A4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf"
$ws.Range("D4").Value = "2016-03-08 21:25:46"
$ws.Range("E4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md"
$ws.Range("F4").Value = "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf"
$ws.Range("G4").Value = "2016-03-08 21:25:06"

# Rebuild hyperlinks: same target URLs, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md", "", "", "ffff442e56c6-4387-47f3-a4e8-312386de752c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ced24756529325343066d6176ecc1533bd35012b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/101f1e70b683caab6d41b2e8adc13f24769d2427/e2e/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3b2d57bce5d1cc161ce3fbca7dc271b578e1ce24/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffff442e56c6-4387-47f3-a4e8-312386de752c.md", "", "", "ffffff60694c53-22b2-4af5-9054-0144e0a92462.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d66aa5de73740e191ff945750899e9a6e616448/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/325f0593420a37db771c4e6075a19e8408a9dad9/e2e/314fba8a-d786-4633-a6c9-f710a8d5fa37.md", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b9b35a0cddba6f6f0f29c2aeb03375210b5dfe6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf", "", "", "314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/e2e/ffffff60694c53-22b2-4af5-9054-0144e0a92462.md", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d66aa5de73740e191ff945750899e9a6e616448/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/325f0593420a37db771c4e6075a19e8408a9dad9/e2e/314fba8a-d786-4633-a6c9-f710a8d5fa37.md", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b9b35a0cddba6f6f0f29c2aeb03375210b5dfe6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/314fba8a-d786-4633-a6c9-f710a8d5fa37.b9059168354fb539a1209ab29727424ed15637ed.de-de.xlf", "", "", "2c46e2e9-de8b-4c6e-8ddc-a579849e99ce.d512c29f63c04efd436e0162f979a8b4733e9556.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c298d950fba8c28142ca911efb99d832db5dc934/.localization-config", "", "", ".localization-config") | Out-Null

